# Auto-generated edit script applying numeric corrections described in the commit diff.
$wb = $excel.ActiveWorkbook

# ALC!row41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 452.9091
$ws.Range("J41").Value = 530.7778
$ws.Range("L41").Value = 530.7778
$ws.Range("N41").Value = -1410.7778

# ALC!row100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2183.75
$ws.Range("I100").Value = 1913.3334
$ws.Range("J100").Value = 2995
$ws.Range("K100").Value = 1913.3334
$ws.Range("L100").Value = 2995
$ws.Range("M100").Value = -1372.3334
$ws.Range("N100").Value = -4077

# ALC!row107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 861.75
$ws.Range("I107").Value = 765.44446
$ws.Range("J107").Value = 985.5714
$ws.Range("K107").Value = 765.44446
$ws.Range("L107").Value = 985.5714
$ws.Range("M107").Value = 1154.55554
$ws.Range("N107").Value = -4825.5714

# ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 5753.5
$ws.Range("I111").Value = 1507
$ws.Range("J111").Value = 10000
$ws.Range("K111").Value = 4521
$ws.Range("L111").Value = 30000
$ws.Range("M111").Value = -1454
$ws.Range("N111").Value = -36134

# ALC!row121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 9001
$ws.Range("J121").Value = 9586.786
$ws.Range("L121").Value = 28760.358
$ws.Range("N121").Value = -32254.358

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 840.0349
$ws.Range("I129").Value = 850
$ws.Range("J129").Value = 839.7976
$ws.Range("K129").Value = 2550
$ws.Range("L129").Value = 2519.3928
$ws.Range("M129").Value = 2450
$ws.Range("N129").Value = -12519.3928

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3960.65
$ws.Range("I132").Value = 4423.9414
$ws.Range("J132").Value = 1335.3334
$ws.Range("K132").Value = 13271.8242
$ws.Range("L132").Value = 4006.0002
$ws.Range("M132").Value = -10741.8242
$ws.Range("N132").Value = -9066.0002

# ALC!row133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 50445
$ws.Range("J133").Value = 50445
$ws.Range("L133").Value = 50445
$ws.Range("N133").Value = -60565

# ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4474.75
$ws.Range("I141").Value = 3000
$ws.Range("K141").Value = 9000
$ws.Range("M141").Value = -3820

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2992
$ws.Range("I32").Value = 2385.0789
$ws.Range("K32").Value = 2385.0789
$ws.Range("M32").Value = -2098.0789

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1939.5333
$ws.Range("I110").Value = 1775
$ws.Range("J110").Value = 2597.6667
$ws.Range("K110").Value = 1775
$ws.Range("L110").Value = 2597.6667
$ws.Range("M110").Value = 270
$ws.Range("N110").Value = -6687.6667

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1433.8889
$ws.Range("I122").Value = 1562.4615
$ws.Range("J122").Value = 1099.6
$ws.Range("K122").Value = 4687.3845
$ws.Range("L122").Value = 3298.8
$ws.Range("M122").Value = -2237.3845
$ws.Range("N122").Value = -8198.799999999999

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2073.8572
$ws.Range("I105").Value = 1937.1111
$ws.Range("J105").Value = 2320
$ws.Range("K105").Value = 1937.1111
$ws.Range("L105").Value = 2320
$ws.Range("M105").Value = -190.1111000000001
$ws.Range("N105").Value = -5814

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2228
$ws.Range("I107").Value = 862.7143
$ws.Range("K107").Value = 862.7143
$ws.Range("M107").Value = 1057.2857

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6605.636
$ws.Range("I134").Value = 6605.636
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 19816.908
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -17281.908
$ws.Range("N134").ClearContents()

# CRP!row105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 12501250
$ws.Range("I105").Value = 17857638
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 17857638
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -17855891
$ws.Range("N105").Value = -6505

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1407.2222
$ws.Range("I134").Value = 1155.7778
$ws.Range("J134").Value = 1658.6666
$ws.Range("K134").Value = 3467.3334
$ws.Range("L134").Value = 4975.9998
$ws.Range("M134").Value = -932.3334000000004
$ws.Range("N134").Value = -10045.9998

# CUL!row11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 110
$ws.Range("I11").Value = 10
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 30
$ws.Range("L11").Value = 480
$ws.Range("M11").Value = 110
$ws.Range("N11").Value = -760

# CUL!row12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 76.28570999999999
$ws.Range("I12").Value = 20
$ws.Range("J12").Value = 98.8
$ws.Range("K12").Value = 60
$ws.Range("L12").Value = 296.4
$ws.Range("M12").Value = 113
$ws.Range("N12").Value = -642.4

# CUL!row70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3292.6875
$ws.Range("I70").Value = 2259.9
$ws.Range("K70").Value = 6779.700000000001
$ws.Range("M70").Value = -6464.700000000001

# CUL!row73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3292.6875
$ws.Range("I73").Value = 2259.9
$ws.Range("K73").Value = 6779.700000000001
$ws.Range("M73").Value = -5687.700000000001

# CUL!row76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4457.778
$ws.Range("J76").Value = 4890
$ws.Range("L76").Value = 14670
$ws.Range("N76").Value = -15436

# CUL!row79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 4457.778
$ws.Range("J79").Value = 4890
$ws.Range("L79").Value = 14670
$ws.Range("N79").Value = -17322

# CUL!row109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1987.5555
$ws.Range("I109").Value = 1129.3334
$ws.Range("J109").Value = 3704
$ws.Range("K109").Value = 3388.0002
$ws.Range("L109").Value = 11112
$ws.Range("M109").Value = -2348.0002
$ws.Range("N109").Value = -13192

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 694.5700000000001
$ws.Range("I131").Value = 302
$ws.Range("J131").Value = 715.23157
$ws.Range("K131").Value = 906
$ws.Range("L131").Value = 2145.69471
$ws.Range("M131").Value = 4134
$ws.Range("N131").Value = -12225.69471

# CUL!row140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1953.1818
$ws.Range("I140").Value = 819.6667
$ws.Range("J140").Value = 3313.4
$ws.Range("K140").Value = 2459.0001
$ws.Range("L140").Value = 9940.200000000001
$ws.Range("M140").Value = 2720.9999
$ws.Range("N140").Value = -20300.2

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4166.6665
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 27320.715
$ws.Range("I132").Value = 3584.4119
$ws.Range("K132").Value = 10753.2357
$ws.Range("M132").Value = -8223.235700000001

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5275
$ws.Range("J22").Value = 5399.5
$ws.Range("L22").Value = 5399.5
$ws.Range("N22").Value = -5989.5

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5275
$ws.Range("J27").Value = 5399.5
$ws.Range("L27").Value = 5399.5
$ws.Range("N27").Value = -5613.5

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 202.71428
$ws.Range("J55").Value = 255.75
$ws.Range("L55").Value = 255.75
$ws.Range("N55").Value = -601.75

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 894479.8
$ws.Range("I122").Value = 1784932.4
$ws.Range("J122").Value = 4027.2727
$ws.Range("K122").Value = 5354797.199999999
$ws.Range("L122").Value = 12081.8181
$ws.Range("M122").Value = -5352347.199999999
$ws.Range("N122").Value = -16981.8181

# WVR!row96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6233.3335
$ws.Range("J96").Value = 7080
$ws.Range("L96").Value = 7080
$ws.Range("N96").Value = -9826

# WVR!row100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1027.0625
$ws.Range("I100").Value = 623
$ws.Range("K100").Value = 1246
$ws.Range("M100").Value = -705

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1379.6389
$ws.Range("I136").Value = 980.1818
$ws.Range("K136").Value = 2940.5454
$ws.Range("M136").Value = -390.5454
